# Commit: "added drp_base to root config"
#
# Adds a new parameter row (drp_base / 0x464000000 / DRP interface base
# address, AXI systems only) to the vu13p_gty_root_config sheet, right
# after the existing mem_base row, and tidies up a leftover formatting
# inconsistency on the row above it (root_config's MGT value cell picks
# up the same style used by the rest of column B). Also updates the saved
# cursor position on that sheet to where it ended up after the edit (A24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix stray style on B18 (root_config's MGT value) so that it uses the
#     same style as the rest of column B's "value" entries (style of B19).
$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial(-4122)

# --- Add new row 23: drp_base parameter, placed right after mem_base (row 22)
#     Start by matching the formatting of the row above (mem_base, row 22).
$ws.Range("A22:C22").Copy()
$ws.Range("A23:C23").PasteSpecial(-4122)

$ws.Range("A23").Value = "drp_base"
$ws.Range("B23").Value = "0x464000000"
$ws.Range("C23").Value = "DRP interface base address, AXI systems only"

# --- Update the saved selection/cursor for this sheet to A24
$ws.Range("A24").Select()
